# Update exception case worksheet with an "Alternate Workflow" column
# and two new ticket rows (AEAREP-645, AEAREP-707), per commit:
# "update exception case with alternate workflow"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two new rows for AEAREP-645 / AEAREP-707 right after AEAREP-471
# (old row 4), shifting the remaining rows down.
$ws.Rows.Item(5).Resize(2).Insert()

# New tickets first (keeps shared-string order: AEAREP-645, AEAREP-707, ...)
$ws.Cells.Item(5, 1).Value = "AEAREP-645"
$ws.Cells.Item(6, 1).Value = "AEAREP-707"

# Header row - rename status_change -> status_change_deviated.
$ws.Range("B1").Value = "status_change_deviated"

$ws.Cells.Item(5, 2).Value = "Approved -> Done"
$ws.Cells.Item(6, 2).Value = "Approved -> Done"

$ws.Cells.Item(5, 3).Value = "Alternate Workflow"
$ws.Cells.Item(6, 3).Value = "Alternate Workflow"

# Fix up the old lowercase "approved -> Done" entry (previously row 4) to use
# the consistent capitalized form already used elsewhere in the column.
$ws.Cells.Item(4, 2).Value = "Approved -> Done"

# Widen column C to fit the new "Alternate Workflow" content (bestFit-style
# width of 46, matching columns A/B which are also sized to fit their text).
$ws.Columns.Item(3).ColumnWidth = 45.1666666666667

# Update selection to match target workbook (C10 selected).
$ws.Range("C10").Select()
